# Append two new rows (68 and 69) of data to the Optical_Power sheet,
# matching the rows already present (A:L are text, M:N are numeric).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        A = "6298"; B = "7/1/2025"; C = "RIVERA INDARTE AV. 1406"; D = "7";
        E = "807877127"; F = "Optical Power"; G = "Pendiente";
        H = "Columna en plantera Se recomienda el retiro riesgo de caida";
        I = "1"; J = "Cambio"; K = "Sin equipos"; L = "Pasante";
        M = -58.450359; N = -34.643582
    },
    @{
        A = "6303"; B = "7/1/2025"; C = "BILBAO, FRANCISCO 2362"; D = "7";
        E = "807877145"; F = "Optical Power"; G = "Pendiente";
        H = "Columna con base corroida oxidada";
        I = "1"; J = "Cambio"; K = "Sin equipos"; L = "Pasante";
        M = -58.459566; N = -34.634615
    }
)

$textCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L")
$rowIndex = 68

foreach ($rowData in $newRows) {
    foreach ($col in $textCols) {
        $cell = $ws.Range("$col$rowIndex")
        # Force text storage (so numeric-/date-looking strings like "6298"
        # or "7/1/2025" are not reinterpreted as a number/date), then drop
        # back to the default "Normal" style so no explicit style index is
        # left on the cell (matching the rest of the data rows).
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$col]
        $cell.Style = "Normal"
    }
    $ws.Range("M$rowIndex").Value = $rowData["M"]
    $ws.Range("N$rowIndex").Value = $rowData["N"]
    $rowIndex++
}
